$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.425.86"
$ws.Range("E2").Value = "  +8.90%  "
$ws.Range("D3").Value = "'1.602.28"
$ws.Range("E3").Value = "  +8.23%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'304.47"
$ws.Range("E5").Value = "  +9.01%  "
$ws.Range("D6").Value = "'0.9921"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").Value = "'0.3683"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "'0.3406"
$ws.Range("E8").Value = "  +10.53%  "
$ws.Range("D9").Value = "'42.73"
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").Value = "'1.141"
$ws.Range("E10").Value = "  +7.04%  "
$ws.Range("D11").Value = "'0.07050"
$ws.Range("E11").Value = "  +5.57%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'5.936"
$ws.Range("E13").Value = "  +7.39%  "
$ws.Range("D14").Value = "'19.66"
$ws.Range("E14").Value = "  +8.64%  "
$ws.Range("D15").Value = "'6.629"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").Value = "'0.00001084"
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").Value = "'1.602.24"
$ws.Range("E17").Value = "  +8.34%  "
$ws.Range("D18").Value = "'0.9919"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "'0.06647"
$ws.Range("E19").Value = "  +11.86%  "
$ws.Range("D20").Value = "'77.97"
$ws.Range("E20").Value = "  +11.82%  "
$ws.Range("D21").Value = "'16.10"
$ws.Range("E21").Value = "  +10.77%  "
$ws.Range("D22").Value = "'6.024"
$ws.Range("E22").Value = "  +9.51%  "
$ws.Range("D23").Value = "'11.80"
$ws.Range("E23").Value = "  +6.67%  "
$ws.Range("D24").Value = "'22.447.20"
$ws.Range("E24").Value = "  +8.77%  "
$ws.Range("E25").Value = "  +6.36%  "
$ws.Range("D26").Value = "'2.509"
$ws.Range("E26").Value = "  +17.50%  "
$ws.Range("D27").Value = "'150.60"
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("D28").Value = "'19.50"
$ws.Range("E28").Value = "  +12.82%  "
$ws.Range("D29").Value = "'1.780.69"
$ws.Range("E29").Value = "  +8.82%  "
$ws.Range("D30").Value = "'120.64"
$ws.Range("E30").Value = "  +5.64%  "
$ws.Range("D31").Value = "'4.198"
$ws.Range("E31").Value = "  +6.78%  "
$ws.Range("D32").Value = "'6.031"
$ws.Range("E32").Value = "  +19.98%  "
$ws.Range("D33").Value = "'0.9570"
$ws.Range("E33").Value = "  +16.23%  "
$ws.Range("D34").Value = "'0.08277"
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("D35").Value = "'1.637"
$ws.Range("E35").Value = "  +6.49%  "
$ws.Range("D36").Value = "'5.288"
$ws.Range("E36").Value = "  +11.82%  "
$ws.Range("D37").Value = "'1.274"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("D38").Value = "'11.84"
$ws.Range("E38").Value = "  +12.91%  "
$ws.Range("D39").Value = "'8.591"
$ws.Range("E39").Value = "  +12.14%  "
$ws.Range("D40").Value = "'0.06107"
$ws.Range("E40").Value = "  +5.57%  "
$ws.Range("D41").Value = "'0.02219"
$ws.Range("E41").Value = "  +8.41%  "
$ws.Range("D42").Value = "'0.2027"
$ws.Range("E42").Value = "  +7.51%  "
$ws.Range("D43").Value = "'0.9914"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("D44").Value = "'0.5910"
$ws.Range("E44").Value = "  +11.38%  "
$ws.Range("E45").Value = "  +9.28%  "
$ws.Range("D46").Value = "'13.26"
$ws.Range("E46").Value = "  +8.62%  "
$ws.Range("D47").Value = "'0.5690"
$ws.Range("E47").Value = "  +9.43%  "
$ws.Range("D48").Value = "'127.01"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("D49").Value = "'1.967"
$ws.Range("E49").Value = "  +8.78%  "
$ws.Range("D50").Value = "'0.06818"
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("D51").Value = "'73.80"
$ws.Range("E51").Value = "  +9.24%  "
